$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 was stored as a text value "76442781"; convert it to a real number.
$ws.Cells.Item(3, 1).Value = 76442781

# Add new row 4 with the new payment record.
# Phone numbers are stored as text, so force text formatting, assign the
# value, then clear the formatting override so no extra style is left on
# the cell (matches the plain/unstyled cells used elsewhere in the sheet).
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "71652621"
$ws.Cells.Item(4, 1).ClearFormats()

$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = "Cash"
$ws.Cells.Item(4, 4).Value = "2025-08-15T09:33:45"
